$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("24 Apr 2020", "0,6", "-10,6", "-22.5,6", "-47.5,6", "-27.5,6", "-17.5,6", "-40,6", "-22.5,4,8,7", "-25,4,8,6"),
    @("28 Apr 2020", "0,6", "-10,6", "-10,6", "-45,3,7,5,5", "-27.5,5,9", "-17.5,4,8,8", "-40,4,7,6", "-22.5,3,7,5,5", "-25,4,5,6"),
    @("7 May 2020", "0,6", "-10,6", "-20,6", "-45,3,8,6,4", "-27.5,6", "-17.5,5,9", "-40,4,8,7", "-22.5,5,6", "-25,3,6,5,7")
)

$startRow = 25
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $rowData[$c]
    }
}

[void]$ws.Range("J27").Select()
